# Auto-generated COM-interop edit script
# Applies per-cell value updates / clears to match the target diff
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# ALC: 27 cell operations
$edits = @(
  @(51, 8, "set", 2994.1177),
  @(51, 9, "set", 2994.1177),
  @(51, 10, "set", 0),
  @(51, 11, "set", 2994.1177),
  @(51, 12, "set", 0),
  @(51, 13, "set", -2510.1177),
  @(51, 14, "clear", $null),
  @(53, 8, "set", 222.1),
  @(53, 10, "set", 186.27272),
  @(53, 12, "set", 186.27272),
  @(53, 14, "set", -1460.27272),
  @(86, 8, "set", 6683.0713),
  @(86, 9, "set", 3415.0908),
  @(86, 11, "set", 3415.0908),
  @(86, 13, "set", -2292.0908),
  @(89, 8, "set", 6683.0713),
  @(89, 9, "set", 3415.0908),
  @(89, 11, "set", 17075.454),
  @(89, 13, "set", -11459.454),
  @(132, 8, "set", 19253.916),
  @(132, 9, "set", 19253.916),
  @(132, 11, "set", 57761.74800000001),
  @(132, 13, "set", -55231.74800000001),
  @(137, 8, "set", 22733170),
  @(137, 10, "set", 12255.223),
  @(137, 12, "set", 36765.669),
  @(137, 14, "set", -41865.669)
)
foreach ($edit in $edits) {
  $r = $edit[0]; $c = $edit[1]; $kind = $edit[2]; $val = $edit[3]
  if ($kind -eq "set") {
    $ws.Cells.Item($r, $c).Value = $val
  } else {
    $ws.Cells.Item($r, $c).ClearContents()
  }
}

$ws = $wb.Worksheets.Item("ARM")
# ARM: 31 cell operations
$edits = @(
  @(32, 8, "set", 1239491.6),
  @(32, 9, "set", 1357259.1),
  @(32, 11, "set", 1357259.1),
  @(32, 13, "set", -1356972.1),
  @(61, 8, "set", 5268020),
  @(61, 9, "set", 4730.25),
  @(61, 11, "set", 4730.25),
  @(61, 13, "set", -4518.25),
  @(74, 8, "set", 1470865.9),
  @(74, 9, "set", 2142174.2),
  @(74, 11, "set", 2142174.2),
  @(74, 13, "set", -2141300.2),
  @(77, 8, "set", 1470865.9),
  @(77, 9, "set", 2142174.2),
  @(77, 11, "set", 10710871),
  @(77, 13, "set", -10706503),
  @(97, 8, "set", 1394.8182),
  @(97, 9, "set", 1214.3),
  @(97, 10, "set", 3200),
  @(97, 11, "set", 1214.3),
  @(97, 12, "set", 3200),
  @(97, 13, "set", -718.3),
  @(97, 14, "set", -4192),
  @(132, 8, "set", 3540.923),
  @(132, 9, "set", 2048.0588),
  @(132, 11, "set", 6144.176399999999),
  @(132, 13, "set", -3614.176399999999),
  @(136, 8, "set", 5268020),
  @(136, 9, "set", 4730.25),
  @(136, 11, "set", 14190.75),
  @(136, 13, "set", -11640.75)
)
foreach ($edit in $edits) {
  $r = $edit[0]; $c = $edit[1]; $kind = $edit[2]; $val = $edit[3]
  if ($kind -eq "set") {
    $ws.Cells.Item($r, $c).Value = $val
  } else {
    $ws.Cells.Item($r, $c).ClearContents()
  }
}

$ws = $wb.Worksheets.Item("BSM")
# BSM: 12 cell operations
$edits = @(
  @(105, 8, "set", 2278.4814),
  @(105, 9, "set", 1973.5294),
  @(105, 11, "set", 1973.5294),
  @(105, 13, "set", -226.5293999999999),
  @(132, 8, "set", 74136.336),
  @(132, 10, "set", 74136.336),
  @(132, 12, "set", 74136.336),
  @(132, 14, "set", -84256.336),
  @(134, 8, "set", 4278682),
  @(134, 9, "set", 2572.8333),
  @(134, 11, "set", 7718.499899999999),
  @(134, 13, "set", -5183.499899999999)
)
foreach ($edit in $edits) {
  $r = $edit[0]; $c = $edit[1]; $kind = $edit[2]; $val = $edit[3]
  if ($kind -eq "set") {
    $ws.Cells.Item($r, $c).Value = $val
  } else {
    $ws.Cells.Item($r, $c).ClearContents()
  }
}

$ws = $wb.Worksheets.Item("CRP")
# CRP: 50 cell operations
$edits = @(
  @(31, 8, "set", 1646262.9),
  @(31, 9, "set", 2507716.8),
  @(31, 10, "set", 1669.1818),
  @(31, 11, "set", 2507716.8),
  @(31, 12, "set", 1669.1818),
  @(31, 13, "set", -2507421.8),
  @(31, 14, "set", -2259.1818),
  @(34, 8, "set", 1646262.9),
  @(34, 9, "set", 2507716.8),
  @(34, 10, "set", 1669.1818),
  @(34, 11, "set", 2507716.8),
  @(34, 12, "set", 1669.1818),
  @(34, 13, "set", -2507514.8),
  @(34, 14, "set", -2073.1818),
  @(58, 8, "set", 7521599),
  @(58, 9, "set", 10420936),
  @(58, 10, "set", 4208070.5),
  @(58, 11, "set", 10420936),
  @(58, 12, "set", 4208070.5),
  @(58, 13, "set", -10420733),
  @(58, 14, "set", -4208476.5),
  @(99, 8, "set", 11281.5),
  @(99, 9, "set", 18294.54),
  @(99, 10, "set", 2993.3635),
  @(99, 11, "set", 18294.54),
  @(99, 12, "set", 2993.3635),
  @(99, 13, "set", -16796.54),
  @(99, 14, "set", -5989.363499999999),
  @(122, 8, "set", 6220.018),
  @(122, 9, "set", 1777.7333),
  @(122, 11, "set", 5333.199900000001),
  @(122, 13, "set", -2883.199900000001),
  @(126, 8, "set", 11281.5),
  @(126, 9, "set", 18294.54),
  @(126, 10, "set", 2993.3635),
  @(126, 11, "set", 54883.62),
  @(126, 12, "set", 8980.0905),
  @(126, 13, "set", -52413.62),
  @(126, 14, "set", -13920.0905),
  @(132, 8, "set", 5166.9585),
  @(132, 10, "set", 6799),
  @(132, 12, "set", 20397),
  @(132, 14, "set", -25457),
  @(136, 8, "set", 7521599),
  @(136, 9, "set", 10420936),
  @(136, 10, "set", 4208070.5),
  @(136, 11, "set", 31262808),
  @(136, 12, "set", 12624211.5),
  @(136, 13, "set", -31260258),
  @(136, 14, "set", -12629311.5)
)
foreach ($edit in $edits) {
  $r = $edit[0]; $c = $edit[1]; $kind = $edit[2]; $val = $edit[3]
  if ($kind -eq "set") {
    $ws.Cells.Item($r, $c).Value = $val
  } else {
    $ws.Cells.Item($r, $c).ClearContents()
  }
}

$ws = $wb.Worksheets.Item("CUL")
# CUL: 32 cell operations
$edits = @(
  @(60, 8, "set", 19999),
  @(60, 9, "set", 19999),
  @(60, 10, "set", 0),
  @(60, 11, "set", 59997),
  @(60, 12, "set", 0),
  @(60, 13, "set", -59746),
  @(60, 14, "clear", $null),
  @(122, 8, "set", 1614827),
  @(122, 10, "set", 2678.4285),
  @(122, 12, "set", 24105.8565),
  @(122, 14, "set", -29005.8565),
  @(123, 8, "set", 13062.5),
  @(123, 9, "set", 7500),
  @(123, 10, "set", 14916.667),
  @(123, 11, "set", 22500),
  @(123, 12, "set", 44750.001),
  @(123, 13, "set", -20050),
  @(123, 14, "set", -49650.001),
  @(124, 8, "set", 24124.223),
  @(124, 9, "set", 35529.5),
  @(124, 10, "set", 15000),
  @(124, 11, "set", 106588.5),
  @(124, 12, "set", 45000),
  @(124, 13, "set", -101678.5),
  @(124, 14, "set", -54820),
  @(126, 8, "set", 14494.833),
  @(126, 9, "set", 0),
  @(126, 10, "set", 14494.833),
  @(126, 11, "set", 0),
  @(126, 12, "set", 43484.499),
  @(126, 13, "clear", $null),
  @(126, 14, "set", -53364.499)
)
foreach ($edit in $edits) {
  $r = $edit[0]; $c = $edit[1]; $kind = $edit[2]; $val = $edit[3]
  if ($kind -eq "set") {
    $ws.Cells.Item($r, $c).Value = $val
  } else {
    $ws.Cells.Item($r, $c).ClearContents()
  }
}

$ws = $wb.Worksheets.Item("GSM")
# GSM: 23 cell operations
$edits = @(
  @(38, 8, "set", 26249.5),
  @(38, 10, "set", 26249.5),
  @(38, 12, "set", 26249.5),
  @(38, 14, "set", -27175.5),
  @(70, 8, "set", 29826.117),
  @(70, 9, "set", 28003.143),
  @(70, 11, "set", 28003.143),
  @(70, 13, "set", -27733.143),
  @(73, 8, "set", 29826.117),
  @(73, 9, "set", 28003.143),
  @(73, 11, "set", 28003.143),
  @(73, 13, "set", -27067.143),
  @(126, 8, "set", 3029.3635),
  @(126, 9, "set", 2669.3333),
  @(126, 11, "set", 8007.999899999999),
  @(126, 13, "set", -5537.999899999999),
  @(132, 8, "set", 11001.667),
  @(132, 9, "set", 12599.7),
  @(132, 10, "set", 3011.5),
  @(132, 11, "set", 37799.10000000001),
  @(132, 12, "set", 9034.5),
  @(132, 13, "set", -35269.10000000001),
  @(132, 14, "set", -14094.5)
)
foreach ($edit in $edits) {
  $r = $edit[0]; $c = $edit[1]; $kind = $edit[2]; $val = $edit[3]
  if ($kind -eq "set") {
    $ws.Cells.Item($r, $c).Value = $val
  } else {
    $ws.Cells.Item($r, $c).ClearContents()
  }
}

$ws = $wb.Worksheets.Item("LTW")
# LTW: 46 cell operations
$edits = @(
  @(40, 8, "set", 4301.6),
  @(40, 9, "set", 3377),
  @(40, 10, "set", 8000),
  @(40, 11, "set", 3377),
  @(40, 12, "set", 8000),
  @(40, 13, "set", -3241),
  @(40, 14, "set", -8272),
  @(80, 8, "set", 44660),
  @(80, 10, "set", 44660),
  @(80, 12, "set", 44660),
  @(80, 14, "set", -46906),
  @(82, 8, "set", 3038),
  @(82, 9, "set", 3057.25),
  @(82, 11, "set", 3057.25),
  @(82, 13, "set", -2696.25),
  @(83, 8, "set", 44660),
  @(83, 10, "set", 44660),
  @(83, 12, "set", 133980),
  @(83, 14, "set", -145212),
  @(85, 8, "set", 3038),
  @(85, 9, "set", 3057.25),
  @(85, 11, "set", 3057.25),
  @(85, 13, "set", -1809.25),
  @(122, 8, "set", 5284.381),
  @(122, 9, "set", 4153.5386),
  @(122, 10, "set", 7122),
  @(122, 11, "set", 12460.6158),
  @(122, 12, "set", 21366),
  @(122, 13, "set", -10010.6158),
  @(122, 14, "set", -26266),
  @(130, 8, "set", 95466.664),
  @(130, 10, "set", 95466.664),
  @(130, 12, "set", 95466.664),
  @(130, 14, "set", -105506.664),
  @(132, 8, "set", 6065251.5),
  @(132, 9, "set", 7411306.5),
  @(132, 11, "set", 22233919.5),
  @(132, 13, "set", -22231389.5),
  @(133, 8, "set", 88913.5),
  @(133, 10, "set", 88913.5),
  @(133, 12, "set", 88913.5),
  @(133, 14, "set", -93973.5),
  @(136, 8, "set", 7355453.5),
  @(136, 9, "set", 4631814),
  @(136, 11, "set", 13895442),
  @(136, 13, "set", -13892892)
)
foreach ($edit in $edits) {
  $r = $edit[0]; $c = $edit[1]; $kind = $edit[2]; $val = $edit[3]
  if ($kind -eq "set") {
    $ws.Cells.Item($r, $c).Value = $val
  } else {
    $ws.Cells.Item($r, $c).ClearContents()
  }
}

$ws = $wb.Worksheets.Item("WVR")
# WVR: 32 cell operations
$edits = @(
  @(74, 8, "set", 7700746.5),
  @(74, 9, "set", 50002110),
  @(74, 10, "set", 9589.272000000001),
  @(74, 11, "set", 50002110),
  @(74, 12, "set", 9589.272000000001),
  @(74, 13, "set", -50001174),
  @(74, 14, "set", -11461.272),
  @(77, 8, "set", 7700746.5),
  @(77, 9, "set", 50002110),
  @(77, 10, "set", 9589.272000000001),
  @(77, 11, "set", 150006330),
  @(77, 12, "set", 28767.816),
  @(77, 13, "set", -150001650),
  @(77, 14, "set", -38127.81600000001),
  @(96, 8, "set", 1184.2106),
  @(96, 9, "set", 1294.0834),
  @(96, 10, "set", 995.8570999999999),
  @(96, 11, "set", 1294.0834),
  @(96, 12, "set", 995.8570999999999),
  @(96, 13, "set", 78.91660000000002),
  @(96, 14, "set", -3741.8571),
  @(122, 8, "set", 24775.82),
  @(122, 9, "set", 2011.775),
  @(122, 11, "set", 6035.325000000001),
  @(122, 13, "set", -3585.325000000001),
  @(136, 8, "set", 3213180.5),
  @(136, 9, "set", 1673714.5),
  @(136, 10, "set", 8931197),
  @(136, 11, "set", 5021143.5),
  @(136, 12, "set", 26793591),
  @(136, 13, "set", -5018593.5),
  @(136, 14, "set", -26798691)
)
foreach ($edit in $edits) {
  $r = $edit[0]; $c = $edit[1]; $kind = $edit[2]; $val = $edit[3]
  if ($kind -eq "set") {
    $ws.Cells.Item($r, $c).Value = $val
  } else {
    $ws.Cells.Item($r, $c).ClearContents()
  }
}

Write-Output "Applied all cell edits"